$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting existing B,C,D to C,D,E
$ws.Columns.Item(2).Insert()

# Give the new header cell the same bold/centered/bordered style as the
# rest of the header row (copy formatting from C1, the old B1 header)
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "segments"

# For each data row, move the segment name from column A to column B,
# and place a numeric 0-based index into column A
for ($r = 2; $r -le 20; $r++) {
    $name = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Column A data cells keep the bold/centered/bordered style (same as before),
# column B (now holding plain text) should not have that style
$ws.Range("B2:B20").Style = "Normal"
